$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Bangkrang Nonthaburi")
$ws.Range("E2").Value = 0.038
$ws.Range("E3").Value = 0.038
$ws.Range("E4").Value = 0.038
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.0254
$ws.Range("O4").Value = 0.00633333333333333
$ws.Range("P4").Value = 0.00633333333333333
$ws.Range("Q4").Value = 0.00633333333333333
$ws.Range("R4").Value = 0.019
$ws.Range("S4").Value = 0.00633333333333333
$ws.Range("T4").Value = 0.00633333333333333
$ws.Range("U4").Value = 0.00633333333333333
$ws.Range("V4").Value = 0.019
$ws.Range("W4").Value = 0.076
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = 0
$ws.Range("E8").Value = 0.0696
$ws.Range("E9").Value = 0.0696
$ws.Range("E10").Value = 0.0696
$ws.Range("M10").Value = 0.0249
$ws.Range("N10").Value = 0.0457
$ws.Range("O10").Value = 0.0116
$ws.Range("P10").Value = 0.0116
$ws.Range("Q10").Value = 0.0116
$ws.Range("R10").Value = 0.0348
$ws.Range("S10").Value = 0.0116
$ws.Range("T10").Value = 0.0116
$ws.Range("U10").Value = 0.0116
$ws.Range("V10").Value = 0.0348
$ws.Range("W10").Value = 0.1392

$ws = $wb.Worksheets.Item("Yueyang China")
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.0263
$ws.Range("O4").Value = 0.00875
$ws.Range("P4").Value = 0.00875
$ws.Range("Q4").Value = 0.00875
$ws.Range("R4").Value = 0.02625
$ws.Range("S4").Value = 0.00875
$ws.Range("T4").Value = 0.00875
$ws.Range("U4").Value = 0.00875
$ws.Range("V4").Value = 0.02625
$ws.Range("W4").Value = 0.105
$ws.Range("E7").Value = 0.0191
$ws.Range("E8").Value = 0.0191
$ws.Range("E9").Value = 0.0191
$ws.Range("M9").Value = 0.0024
$ws.Range("N9").Value = 0.0048
$ws.Range("O9").Value = 0.00318333333333333
$ws.Range("P9").Value = 0.00318333333333333
$ws.Range("Q9").Value = 0.00318333333333333
$ws.Range("R9").Value = 0.00955
$ws.Range("S9").Value = 0.00318333333333333
$ws.Range("T9").Value = 0.00318333333333333
$ws.Range("U9").Value = 0.00318333333333333
$ws.Range("V9").Value = 0.00955
$ws.Range("W9").Value = 0.0382

$ws = $wb.Worksheets.Item("Changzhou Epc China")
$ws.Range("E2").Value = 0.0693
$ws.Range("E3").Value = 0.0693
$ws.Range("E4").Value = 0.0693
$ws.Range("M4").Value = 0.0357
$ws.Range("N4").Value = 0.0524
$ws.Range("O4").Value = 0.01155
$ws.Range("P4").Value = 0.01155
$ws.Range("Q4").Value = 0.01155
$ws.Range("R4").Value = 0.03465
$ws.Range("S4").Value = 0.01155
$ws.Range("T4").Value = 0.01155
$ws.Range("U4").Value = 0.01155
$ws.Range("V4").Value = 0.03465
$ws.Range("W4").Value = 0.1386
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("E8").Value = 0.0522
$ws.Range("E9").Value = 0.0522
$ws.Range("E10").Value = 0.0522
$ws.Range("M10").Value = 0.0035
$ws.Range("N10").Value = 0.0317
$ws.Range("O10").Value = 0.0087
$ws.Range("P10").Value = 0.0087
$ws.Range("Q10").Value = 0.0087
$ws.Range("R10").Value = 0.0261
$ws.Range("S10").Value = 0.0087
$ws.Range("T10").Value = 0.0087
$ws.Range("U10").Value = 0.0087
$ws.Range("V10").Value = 0.0261
$ws.Range("W10").Value = 0.1044

$ws = $wb.Worksheets.Item("Jiaxing China")
$ws.Range("E2").Value = 0.0478
$ws.Range("E3").Value = 0.0478
$ws.Range("E4").Value = 0.0478
$ws.Range("M4").Value = 0.0488
$ws.Range("N4").Value = 0.0476
$ws.Range("O4").Value = 0.00796666666666667
$ws.Range("P4").Value = 0.00796666666666667
$ws.Range("Q4").Value = 0.00796666666666667
$ws.Range("R4").Value = 0.0239
$ws.Range("S4").Value = 0.00796666666666667
$ws.Range("T4").Value = 0.00796666666666667
$ws.Range("U4").Value = 0.00796666666666667
$ws.Range("V4").Value = 0.0239
$ws.Range("W4").Value = 0.0956
$ws.Range("E7").Value = 0.0243
$ws.Range("E8").Value = 0.0243
$ws.Range("E9").Value = 0.0243
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0.0157
$ws.Range("O9").Value = 0.00405
$ws.Range("P9").Value = 0.00405
$ws.Range("Q9").Value = 0.00405
$ws.Range("R9").Value = 0.01215
$ws.Range("S9").Value = 0.00405
$ws.Range("T9").Value = 0.00405
$ws.Range("U9").Value = 0.00405
$ws.Range("V9").Value = 0.01215
$ws.Range("W9").Value = 0.0486

$ws = $wb.Worksheets.Item("Suzhou China")
$ws.Range("E2").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("E5").Value = 0.054
$ws.Range("E6").Value = 0.054
$ws.Range("E7").Value = 0.054
$ws.Range("M7").Value = 0.0066
$ws.Range("N7").Value = 0.0302
$ws.Range("O7").Value = 0.009
$ws.Range("P7").Value = 0.009
$ws.Range("Q7").Value = 0.009
$ws.Range("R7").Value = 0.027
$ws.Range("S7").Value = 0.009
$ws.Range("T7").Value = 0.009
$ws.Range("U7").Value = 0.009
$ws.Range("V7").Value = 0.027
$ws.Range("W7").Value = 0.108
